$d = $word.ActiveDocument

# 1) Merge "2.7 " + "[draw a diagram" + " with labels]" into a single run's text.
$d.Content.Find.Execute("2.7 [draw a diagram with labels]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2.7 [draw a diagram with labels]", 2)

# 2) Merge ">" + "500 Ohm measurement: " + "x" + " Ohms" into a single run's text.
$d.Content.Find.Execute(">500 Ohm measurement: x Ohms", $true, $false, $false, $false, $false,
                         $true, 1, $false, ">500 Ohm measurement: x Ohms", 2)

# 3) Merge ">" + "Resistance of all resistors in series: " + "x" + " Ohms" into a single run's text.
$d.Content.Find.Execute(">Resistance of all resistors in series: x Ohms", $true, $false, $false, $false, $false,
                         $true, 1, $false, ">Resistance of all resistors in series: x Ohms", 2)

# 4) "Arduino - Blink a LED" -> "ESP32" + " - Blink a LED" (bold kept on both runs).
$d.Content.Find.Execute("Arduino", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ESP32", 2)
